$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row inserts in this engine do not carry Hyperlink range refs along with the shifted cells,
# so clear all existing hyperlinks now and rebuild them after every cell value is final.
$ws.Hyperlinks.Delete()

# Insert 3 new listing rows (matches the 3 new rows appended/spliced into the feed):
#  - one brand-new listing lands at (final) row 7, pushing the old row7..row13 down by one
#  - two more brand-new listings land at (final) rows 14-15, pushing the rest down further
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

# Write every cell of the refreshed scrape (header untouched) so values/order are exactly right
# Row 2: '【AI活用】データ分析Webサービス開発パートナ'
$ws.Cells.Item(2, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(2, 2).Value2 = '【AI活用】データ分析Webサービス開発パートナー募集'
$ws.Cells.Item(2, 3).Value2 = 'システム開発'
$ws.Cells.Item(2, 4).Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(2, 5).Value2 = '期限情報なし'
$ws.Cells.Item(2, 6).Value2 = 'https://www.lancers.jp/work/detail/5399092'
$ws.Cells.Item(2, 7).Value2 = 368
$ws.Cells.Item(2, 8).Value2 = '🔥AI,Ai ◆開発'

# Row 3: 'あなたAIクローン構築パートナー募集・モデル制作'
$ws.Cells.Item(3, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(3, 2).Value2 = 'あなたAIクローン構築パートナー募集・モデル制作&新規依頼'
$ws.Cells.Item(3, 3).Value2 = 'システム開発'
$ws.Cells.Item(3, 4).Value2 = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(3, 5).Value2 = '期限情報なし'
$ws.Cells.Item(3, 6).Value2 = 'https://www.lancers.jp/work/detail/5399534'
$ws.Cells.Item(3, 7).Value2 = 303
$ws.Cells.Item(3, 8).Value2 = '🔥AI,Ai'

# Row 4: 'Excel・Accessベースの改修や追加、Py'
$ws.Cells.Item(4, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(4, 2).Value2 = 'Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集'
$ws.Cells.Item(4, 3).Value2 = 'システム開発'
$ws.Cells.Item(4, 4).Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(4, 5).Value2 = '期限情報なし'
$ws.Cells.Item(4, 6).Value2 = 'https://www.lancers.jp/work/detail/5399398'
$ws.Cells.Item(4, 7).Value2 = 298
$ws.Cells.Item(4, 8).Value2 = '🔥Python ◆開発,スクレイピング'

# Row 5: '【急募】カスタマー向けFAQチャットbotの開発'
$ws.Cells.Item(5, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(5, 2).Value2 = '【急募】カスタマー向けFAQチャットbotの開発依頼'
$ws.Cells.Item(5, 3).Value2 = 'システム開発'
$ws.Cells.Item(5, 4).Value2 = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(5, 5).Value2 = '期限情報なし'
$ws.Cells.Item(5, 6).Value2 = 'https://www.lancers.jp/work/detail/5399558'
$ws.Cells.Item(5, 7).Value2 = 180
$ws.Cells.Item(5, 8).Value2 = '★bot ◆開発'

# Row 6: '既存Excelをベースにした短期計画書管理のWe'
$ws.Cells.Item(6, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(6, 2).Value2 = '既存Excelをベースにした短期計画書管理のWebシステム開発'
$ws.Cells.Item(6, 3).Value2 = 'システム開発'
$ws.Cells.Item(6, 4).Value2 = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(6, 5).Value2 = '期限情報なし'
$ws.Cells.Item(6, 6).Value2 = 'https://www.lancers.jp/work/detail/5399602'
$ws.Cells.Item(6, 7).Value2 = 153
$ws.Cells.Item(6, 8).Value2 = '◆開発,システム開発 ◇管理'

# Row 7: '【Flutter+Firebase】社内ポータル'
$ws.Cells.Item(7, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(7, 2).Value2 = '【Flutter+Firebase】社内ポータルアプリ開発のパートナー募集'
$ws.Cells.Item(7, 3).Value2 = 'システム開発'
$ws.Cells.Item(7, 4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(7, 5).Value2 = '期限情報なし'
$ws.Cells.Item(7, 6).Value2 = 'https://www.lancers.jp/work/detail/5399824'
$ws.Cells.Item(7, 7).Value2 = 100
$ws.Cells.Item(7, 8).Value2 = '◆開発 ◇アプリ'

# Row 8: '【RPA/ブラウザ操作自動化】Webフォーム大量'
$ws.Cells.Item(8, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(8, 2).Value2 = '【RPA/ブラウザ操作自動化】Webフォーム大量登録の自動化(継続依頼あり)'
$ws.Cells.Item(8, 3).Value2 = 'システム開発'
$ws.Cells.Item(8, 4).Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(8, 5).Value2 = '期限情報なし'
$ws.Cells.Item(8, 6).Value2 = 'https://www.lancers.jp/work/detail/5399631'
$ws.Cells.Item(8, 7).Value2 = 88
$ws.Cells.Item(8, 8).Value2 = '◆自動化'

# Row 9: 'Googleフォーム × スプレッドシート × '
$ws.Cells.Item(9, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(9, 2).Value2 = 'Googleフォーム × スプレッドシート × GAS 自動化(ストレスチェック診断/台帳保存あり)'
$ws.Cells.Item(9, 3).Value2 = 'システム開発'
$ws.Cells.Item(9, 4).Value2 = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(9, 5).Value2 = '期限情報なし'
$ws.Cells.Item(9, 6).Value2 = 'https://www.lancers.jp/work/detail/5399200'
$ws.Cells.Item(9, 7).Value2 = 88
$ws.Cells.Item(9, 8).Value2 = '◆自動化'

# Row 10: '急募 【急募】Excelで株の保有リストを自動化'
$ws.Cells.Item(10, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(10, 2).Value2 = '急募 【急募】Excelで株の保有リストを自動化したいので制作してくださる方募集!'
$ws.Cells.Item(10, 3).Value2 = 'システム開発'
$ws.Cells.Item(10, 4).Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(10, 5).Value2 = '期限情報なし'
$ws.Cells.Item(10, 6).Value2 = 'https://www.lancers.jp/work/detail/5399727'
$ws.Cells.Item(10, 7).Value2 = 83
$ws.Cells.Item(10, 8).Value2 = '◆自動化'

# Row 11: '【急募】住宅展示場マッチング診断サービスのMVP'
$ws.Cells.Item(11, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(11, 2).Value2 = '【急募】住宅展示場マッチング診断サービスのMVP開発依頼'
$ws.Cells.Item(11, 3).Value2 = 'システム開発'
$ws.Cells.Item(11, 4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(11, 5).Value2 = '期限情報なし'
$ws.Cells.Item(11, 6).Value2 = 'https://www.lancers.jp/work/detail/5399759'
$ws.Cells.Item(11, 7).Value2 = 75
$ws.Cells.Item(11, 8).Value2 = '◆開発'

# Row 12: '完全在宅GASエンジニア募集/課題テストからご依'
$ws.Cells.Item(12, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(12, 2).Value2 = '完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします'
$ws.Cells.Item(12, 3).Value2 = 'システム開発'
$ws.Cells.Item(12, 4).Value2 = '~ 5,000 円 / 固定'
$ws.Cells.Item(12, 5).Value2 = '期限情報なし'
$ws.Cells.Item(12, 6).Value2 = 'https://www.lancers.jp/work/detail/5399071'
$ws.Cells.Item(12, 7).Value2 = 70
$ws.Cells.Item(12, 8).Value2 = '◆効率化'

# Row 13: '【ペットのアバター化】Pawsitiveプロトタ'
$ws.Cells.Item(13, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(13, 2).Value2 = '【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼'
$ws.Cells.Item(13, 3).Value2 = 'システム開発'
$ws.Cells.Item(13, 4).Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(13, 5).Value2 = '期限情報なし'
$ws.Cells.Item(13, 6).Value2 = 'https://www.lancers.jp/work/detail/5399313'
$ws.Cells.Item(13, 7).Value2 = 68
$ws.Cells.Item(13, 8).Value2 = '◆開発'

# Row 14: '【急募】PHP・Lalavelでの既存プログラム'
$ws.Cells.Item(14, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(14, 2).Value2 = '【急募】PHP・Lalavelでの既存プログラム改修依頼'
$ws.Cells.Item(14, 3).Value2 = 'システム開発'
$ws.Cells.Item(14, 4).Value2 = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(14, 5).Value2 = '期限情報なし'
$ws.Cells.Item(14, 6).Value2 = 'https://www.lancers.jp/work/detail/5396563'
$ws.Cells.Item(14, 7).Value2 = 33
$ws.Cells.Item(14, 8).Value2 = '○PHP'

# Row 15: '【SES案件多数/リモート可】フロントエンドエン'
$ws.Cells.Item(15, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(15, 2).Value2 = '【SES案件多数/リモート可】フロントエンドエンジニア募集(HTML/CSS〜モダンFWまで歓迎)'
$ws.Cells.Item(15, 3).Value2 = 'システム開発'
$ws.Cells.Item(15, 4).Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(15, 5).Value2 = '期限情報なし'
$ws.Cells.Item(15, 6).Value2 = 'https://www.lancers.jp/work/detail/5399721'
$ws.Cells.Item(15, 7).Value2 = 25

# Row 16: '〖リモート可〗Delphiエンジニア募集'
$ws.Cells.Item(16, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(16, 2).Value2 = '〖リモート可〗Delphiエンジニア募集'
$ws.Cells.Item(16, 3).Value2 = 'システム開発'
$ws.Cells.Item(16, 4).Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(16, 5).Value2 = '期限情報なし'
$ws.Cells.Item(16, 6).Value2 = 'https://www.lancers.jp/work/detail/5341051'
$ws.Cells.Item(16, 7).Value2 = 25

# Row 17: '【急募】フロントエンドエンジニア募集!(ややWE'
$ws.Cells.Item(17, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(17, 2).Value2 = '【急募】フロントエンドエンジニア募集!(ややWEBコーダー寄り)'
$ws.Cells.Item(17, 3).Value2 = 'システム開発'
$ws.Cells.Item(17, 4).Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(17, 5).Value2 = '期限情報なし'
$ws.Cells.Item(17, 6).Value2 = 'https://www.lancers.jp/work/detail/5399545'
$ws.Cells.Item(17, 7).Value2 = 25

# Row 18: '【相談から実装まで伴走できる方歓迎】介護・福祉×'
$ws.Cells.Item(18, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(18, 2).Value2 = '【相談から実装まで伴走できる方歓迎】介護・福祉×テクノロジー事例収集の仕組みづくり'
$ws.Cells.Item(18, 3).Value2 = 'システム開発'
$ws.Cells.Item(18, 4).Value2 = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(18, 5).Value2 = '期限情報なし'
$ws.Cells.Item(18, 6).Value2 = 'https://www.lancers.jp/work/detail/5398932'
$ws.Cells.Item(18, 7).Value2 = 18

# Row 19: '限定公開 PR 限定公開の仕事'
$ws.Cells.Item(19, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(19, 2).Value2 = '限定公開 PR 限定公開の仕事'
$ws.Cells.Item(19, 3).Value2 = 'システム開発'
$ws.Cells.Item(19, 4).Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(19, 5).Value2 = '期限情報なし'
$ws.Cells.Item(19, 6).Value2 = 'https://www.lancers.jp/work/detail/5399347'
$ws.Cells.Item(19, 7).Value2 = 13

# Row 20: 'Android kotlin 非同期処理の呼び方'
$ws.Cells.Item(20, 1).Value2 = '2025-09-24 18:24:58'
$ws.Cells.Item(20, 2).Value2 = 'Android kotlin 非同期処理の呼び方'
$ws.Cells.Item(20, 3).Value2 = 'システム開発'
$ws.Cells.Item(20, 4).Value2 = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(20, 5).Value2 = '期限情報なし'
$ws.Cells.Item(20, 6).Value2 = 'https://www.lancers.jp/work/detail/5399765'
$ws.Cells.Item(20, 7).Value2 = 10

# Rebuild hyperlinks on column F now that every row/value is in its final place
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5399092") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5399534") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5399398") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5399558") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5399602") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5399824") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5399631") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5399200") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5399727") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5399759") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5399071") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5399313") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5396563") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5399721") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5341051") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5399545") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5398932") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5399347") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5399765") | Out-Null

Write-Output "done"
